# Move repeated GUILayout logic to separate file
#
# translations.xlsx changes:
#  - KeyLimiter sheet: add Korean (column C) translations for the key
#    viewer / key-view-color rows that previously only had KEY + ENGLISH.
#  - PlanetColor sheet: collapse the per-channel BODY_R/BODY_G/BODY_B and
#    TAIL_R/TAIL_G/TAIL_B (plus BODY_HEX/TAIL_HEX) rows into a single
#    BODY and TAIL row each (now a single color-picker value instead of
#    separate R/G/B/Hex fields).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# KeyLimiter: fill in the missing KOREAN column for the key-viewer rows.
# ---------------------------------------------------------------------
$keyLimiter = $wb.Worksheets.Item("KeyLimiter")

$keyLimiter.Range("C8").Value  = "등록된 키들의 키뷰어 보이기"
$keyLimiter.Range("C9").Value  = "크기:"
$keyLimiter.Range("C10").Value = "가로 위치:"
$keyLimiter.Range("C11").Value = "세로 위치:"
$keyLimiter.Range("C12").Value = "누른 키의 테두리 색상:"
$keyLimiter.Range("C13").Value = "누르지 않은 키의 테두리 색상:"
$keyLimiter.Range("C14").Value = "누른 키의 배경 색상:"
$keyLimiter.Range("C15").Value = "누르지 않은 키의 배경 색상:"
$keyLimiter.Range("C16").Value = "누른 키의 텍스트 색상:"
$keyLimiter.Range("C17").Value = "누르지 않은 키의 텍스트 색상:"

# ---------------------------------------------------------------------
# PlanetColor: replace the six R/G/B rows (for Body and Tail) plus the
# Body Hex / Tail Hex rows with a single Body / Tail row each.
# ---------------------------------------------------------------------
$planetColor = $wb.Worksheets.Item("PlanetColor")

# Wipe out the old rows 6-13 (BODY_R, BODY_G, BODY_B, TAIL_R, TAIL_G,
# TAIL_B, BODY_HEX, TAIL_HEX) first.
$planetColor.Range("A6:D13").ClearContents()

# Row 6: BODY
$planetColor.Range("A6").Value = "BODY"
$planetColor.Range("B6").Value = "Body:"
$planetColor.Range("C6").Value = "행성:"
$planetColor.Range("D6").Value = "Cuerpo:"

# Row 7: TAIL
$planetColor.Range("A7").Value = "TAIL"
$planetColor.Range("B7").Value = "Tail:"
$planetColor.Range("C7").Value = "꼬리 파티클:"
$planetColor.Range("D7").Value = "Cola:"
